$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from the row above (row 7) down into the previously-blank
# row 8 so the new data row matches the other populated rows' styling.
$ws.Range("A7:O7").Copy()
$ws.Range("A8:O8").PasteSpecial(-4122)

# "> 160/80" -> "> 160/100" for the BP threshold used by the
# ukb51139_subset.csv rows (E3:E7).
$ws.Range("E3:E7").Value = "> 160/100"

# Subset Dimensions typo fix: "2802 x 1081" -> "2801 x 1081" (B6).
$ws.Range("B6").Value = "2801 x 1081"

# Populate the new results row (row 8)
$ws.Range("A8").Value = "ukb51139_subset.csv"
$ws.Range("B8").Value = "2801 x 1081"
$ws.Range("C8").Value = "all"
$ws.Range("D8").Value = "no event"
$ws.Range("E8").Value = "> 140/80"
$ws.Range("F8").Value = "zscore"
$ws.Range("G8").Value = "median"
$ws.Range("H8").Value = "age, sex"
$ws.Range("I8").Value = 50
$ws.Range("K8").Value = 440
$ws.Range("L8").Value = " -378 & -40.7"
$ws.Range("M8").Value = "8.6 & 7.1"
$ws.Range("N8").Value = 17
$ws.Range("O8").Value = 47.5
